# Update the "想去人数" (interested-count) column F for the refreshed
# scrape snapshot. The same ten events live on both the "展览" sheet and
# the aggregate "全部类型" sheet, so both need the identical bump.

$wb = $excel.ActiveWorkbook

$updates = @{
    4  = 10597
    6  = 963
    7  = 104
    8  = 1306
    9  = 8209
    11 = 461
    15 = 3251
    18 = 742
    21 = 283
    23 = 1699
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
